$d = $word.ActiveDocument

# 1. Remove the entire "Meta description" paragraph that currently follows
#    the "Play Barbarian Gold Slot Game Online for Free" heading.
$metaPara = $d.Paragraphs.Item(2)
if ($metaPara.Range.Text -like "Meta description*") {
    $metaPara.Range.Delete()
} else {
    Write-Host "WARNING: paragraph 2 did not look like the Meta description paragraph:" $metaPara.Range.Text
}

# 2. Insert a new bold paragraph "Play Barbarian Gold Slot Game Online for Free"
#    right before the final ("Prompt: ...") paragraph.
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$lastPara.Range.InsertParagraphBefore()

$newPara = $d.Paragraphs.Item($count)
$newParaXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Barbarian Gold Slot Game Online for Free</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$newPara.Range.InsertXML($newParaXml)

# 3. Replace the text of the final paragraph: drop the "Prompt: ..." image-prompt
#    copy and put the (former) meta-description sentence there instead.
$oldText = 'Prompt: Create a feature image for "Barbarian Gold," an online slot game by Iron Dog Studios, featuring a happy Maya warrior with glasses in a cartoon style. The image should reflect the adventurous and playful nature of the game and showcase the fierce protagonist, a barbarian warrior with an axe, sword, and spiked ball, as he travels through the game''s various levels. The image should include a backdrop of a medieval castle and possibly include some of the game''s symbols, such as the lizard enemy and the shield with the crossed swords. In addition, the image should prominently feature the happy Maya warrior with glasses, who acts as a playful nod to the game''s fun and vibrant atmosphere. The Maya warrior should be shown in a confident and triumphant pose, perhaps holding a bag of gold or in the midst of a winning spin. Overall, the feature image should capture the excitement and thrill of playing "Barbarian Gold" and entice players to join in on the adventure.'
$newText = 'Read our review of Barbarian Gold and play for free. Join the warrior in his adventure to defeat the lizard enemy and steal the treasure.'

$d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
